$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in B2:C5 per the diff
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 12

$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 10

$ws.Range("B4").Value = 0.8

$ws.Range("C5").Value = 25

# Update the active selection to C3 (single cell) as in the diff
$ws.Range("C3").Select()
